$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-06-23 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-24 Tuesday", 2)

# Update the division problems table (5 content rows at table rows 1,5,9,13,17; 5 columns)
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="52÷9="},
    @{Row=1;  Col=2; Text="79÷8="},
    @{Row=1;  Col=3; Text="61÷3="},
    @{Row=1;  Col=4; Text="25÷2="},
    @{Row=1;  Col=5; Text="35÷4="},

    @{Row=5;  Col=1; Text="13÷9="},
    @{Row=5;  Col=2; Text="14÷2="},
    @{Row=5;  Col=3; Text="83÷8="},
    @{Row=5;  Col=4; Text="52÷4="},
    @{Row=5;  Col=5; Text="55÷9="},

    @{Row=9;  Col=1; Text="96÷3="},
    @{Row=9;  Col=2; Text="56÷4="},
    @{Row=9;  Col=3; Text="34÷3="},
    @{Row=9;  Col=4; Text="66÷3="},
    @{Row=9;  Col=5; Text="87÷4="},

    @{Row=13; Col=1; Text="98÷2="},
    @{Row=13; Col=2; Text="96÷3="},
    @{Row=13; Col=3; Text="50÷7="},
    @{Row=13; Col=4; Text="44÷5="},
    @{Row=13; Col=5; Text="36÷9="},

    @{Row=17; Col=1; Text="50÷9="},
    @{Row=17; Col=2; Text="62÷2="},
    @{Row=17; Col=3; Text="10÷3="},
    @{Row=17; Col=4; Text="55÷9="},
    @{Row=17; Col=5; Text="71÷5="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Output "done"
